# Updating the slope and constant for RVR based on the assumption that MCFP is accurate.
#
# Semantic change:
#   - H3:H8 previously computed a feedback formula off the G13:G19 block
#     (=G13, =G14, ...). MCFP is now assumed accurate, so these become
#     plain literal zeros instead of formulas.
#   - D13 previously held a hard-coded literal (7); it now uses the same
#     formula pattern already used by D14:D19 -> (B13-C13)*0.0047.
#
# Everything downstream (E/G/K columns, shared-formula groupings, etc.)
# recalculates naturally from these two edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H3:H8 -> literal 0 (formulas removed)
$ws.Range("H3:H8").Value = 0

# D13 -> formula matching D14:D19's pattern
$ws.Range("D13").Formula = "=(B13-C13)*0.0047"

# Restore the cursor/selection position to match the saved workbook state.
$ws.Range("I19").Select()
